$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.0000824793241918087
$ws.Range("E2").Value = 0.01425322098657489
$ws.Range("G2").Value = 0.001233034767210484
$ws.Range("H2").Value = 0.003044792916625738
$ws.Range("I2").Value = 0.003799600526690483
$ws.Range("J2").Value = 0.004561311565339565
$ws.Range("K2").Value = 0.0004860507324337959
$ws.Range("D3").Value = 0.001491243951022625
$ws.Range("E3").Value = 0.02403086423873901
$ws.Range("G3").Value = 0.00196731323376298
$ws.Range("H3").Value = 0.00533769279718399
$ws.Range("I3").Value = 0.004900505300611258
$ws.Range("J3").Value = 0.009090162348002195
$ws.Range("K3").Value = 0.0007988982833921909
$ws.Range("D4").Value = 0.001304802019149065
$ws.Range("E4").Value = 0.02408290188759565
$ws.Range("G4").Value = 0.001909180544316769
$ws.Range("H4").Value = 0.005269759800285101
$ws.Range("I4").Value = 0.005247276742011309
$ws.Range("J4").Value = 0.008864318020641804
$ws.Range("K4").Value = 0.000744241289794445
$ws.Range("D5").Value = 0.0001092320308089256
$ws.Range("E5").Value = 0.01391712389886379
$ws.Range("G5").Value = 0.001341833733022213
$ws.Range("H5").Value = 0.002403884660452604
$ws.Range("I5").Value = 0.003680076450109482
$ws.Range("J5").Value = 0.004938533529639244
$ws.Range("K5").Value = 0.0004652114585042
$ws.Range("D6").Value = 0.001381068956106901
$ws.Range("E6").Value = 0.03841330064460635
$ws.Range("G6").Value = 0.002221256494522095
$ws.Range("H6").Value = 0.005978557281196117
$ws.Range("I6").Value = 0.02072853874415159
$ws.Range("J6").Value = 0.006161280442029238
$ws.Range("K6").Value = 0.001004157587885857
$ws.Range("D8").Value = 0.0000824793241918087
$ws.Range("E8").Value = 0.01425322098657489
$ws.Range("G8").Value = 0.001233034767210484
$ws.Range("H8").Value = 0.003044792916625738
$ws.Range("I8").Value = 0.003799600526690483
$ws.Range("J8").Value = 0.004561311565339565
$ws.Range("K8").Value = 0.0004860507324337959
$ws.Range("D9").Value = 0.001491243951022625
$ws.Range("E9").Value = 0.02403086423873901
$ws.Range("G9").Value = 0.00196731323376298
$ws.Range("H9").Value = 0.00533769279718399
$ws.Range("I9").Value = 0.004900505300611258
$ws.Range("J9").Value = 0.009090162348002195
$ws.Range("K9").Value = 0.0007988982833921909
$ws.Range("D10").Value = 0.001304802019149065
$ws.Range("E10").Value = 0.02408290188759565
$ws.Range("G10").Value = 0.001909180544316769
$ws.Range("H10").Value = 0.005269759800285101
$ws.Range("I10").Value = 0.005247276742011309
$ws.Range("J10").Value = 0.008864318020641804
$ws.Range("K10").Value = 0.000744241289794445
$ws.Range("D11").Value = 0.0001092320308089256
$ws.Range("E11").Value = 0.01391712389886379
$ws.Range("G11").Value = 0.001341833733022213
$ws.Range("H11").Value = 0.002403884660452604
$ws.Range("I11").Value = 0.003680076450109482
$ws.Range("J11").Value = 0.004938533529639244
$ws.Range("K11").Value = 0.0004652114585042
$ws.Range("D12").Value = 0.001381068956106901
$ws.Range("E12").Value = 0.03841330064460635
$ws.Range("G12").Value = 0.002221256494522095
$ws.Range("H12").Value = 0.005978557281196117
$ws.Range("I12").Value = 0.02072853874415159
$ws.Range("J12").Value = 0.006161280442029238
$ws.Range("K12").Value = 0.001004157587885857
